$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A78").Value = 43811
$ws.Range("B78").Value = 2207.7382817953599
$ws.Range("C78").Value = 2207.0300000000002
$ws.Range("D78").Formula = "=100*(B78-C78)/C78"
$ws.Range("E78").Value = 169
$ws.Range("F78").Value = "New CRM opened 12/11/2020"
Write-Host "done"
